$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 0, 8.974608811992548)
    3 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 0, 5.586269137925634)
    4 = @(0.04271373187048222, 0.306821227259698, 0.1494219747398047, 10.19245300693656, 1, 10.69140994080654)
    5 = @(0.04271373187048222, 0.306821227259698, 3.537761648806719, 0.4942365360607697, 1, 4.381533143997669)
    6 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 6.189590430959694)
    7 = @(0.2917716402565462, 0.04071648406533734, 3.537761648806719, 10.19245300693656, 0, 14.06270278006516)
    8 = @(0.6606524410359556, 0.00006240767534437808, 3.537761648806719, 0.4942365360607697, 1, 4.692713033578789)
    9 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 6.189590430959694)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
